# ADD results from server
# Updates numeric result values on each yearly sheet (2025, 2030, 2035, 2040, 2045, 2050)
# to reflect the latest server-computed figures.

$wb = $excel.ActiveWorkbook

$updates = @(
    ,@("2025", "N2", 7169.226093134127)
    ,@("2025", "O2", 6984.121280850342)

    ,@("2030", "B2", 5707.815717280662)
    ,@("2030", "I2", 44492.05901988943)
    ,@("2030", "L2", 66334.06707325629)
    ,@("2030", "M2", 21991.42050229464)
    ,@("2030", "O2", 12089.30731174489)

    ,@("2035", "A2", 2927.360317916481)
    ,@("2035", "B2", 7940.887964949257)
    ,@("2035", "E2", 67179.99183625776)
    ,@("2035", "I2", 59530.75343380851)
    ,@("2035", "L2", 66334.06707325629)
    ,@("2035", "M2", 25547.11936466757)
    ,@("2035", "N2", 15148.30423606118)
    ,@("2035", "O2", 14771.7086506748)

    ,@("2040", "A2", 2927.360317916481)
    ,@("2040", "B2", 7940.887964949257)
    ,@("2040", "E2", 67179.99183625776)
    ,@("2040", "I2", 59530.75343380851)
    ,@("2040", "L2", 66334.06707325629)
    ,@("2040", "M2", 25547.11936466757)
    ,@("2040", "N2", 15255.98985290295)
    ,@("2040", "O2", 14771.7086506748)

    ,@("2045", "A2", 6352.985609279765)
    ,@("2045", "B2", 7940.887964949257)
    ,@("2045", "E2", 67179.99183625776)
    ,@("2045", "I2", 59530.75343380851)
    ,@("2045", "L2", 66334.06707325629)
    ,@("2045", "M2", 25547.11936466757)
    ,@("2045", "N2", 15803.97116121355)
    ,@("2045", "O2", 17114.26384084568)

    ,@("2050", "A2", 6352.985609279765)
    ,@("2050", "B2", 7940.887964949257)
    ,@("2050", "E2", 67179.99183625776)
    ,@("2050", "I2", 59530.75343380851)
    ,@("2050", "L2", 66334.06707325629)
    ,@("2050", "M2", 25547.11936466757)
    ,@("2050", "N2", 15803.97116121355)
    ,@("2050", "O2", 17114.26384084568)
)

foreach ($item in $updates) {
    $sheetName = $item[0]
    $cellRef = $item[1]
    $val = $item[2]
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Range($cellRef).Value = $val
}
